# Weekly update: a new week's worth of price data for "Acelga" was published
# at the top of this product's data block. Insert two new rows right above
# the existing first data row of that block (row 476) and shift everything
# else down by two rows, then populate the two new rows with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 476 (existing rows 476.. shift down to 478..)
$ws.Range("A476:A477").EntireRow.Insert()

# New row 476: Acelga, Primera, week of 2022-02-18
$ws.Range("A476").Value = 9
$ws.Range("B476").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C476").Value = "Metropolitana"
$ws.Range("D476").Value = 44610
$ws.Range("E476").Value = 13
$ws.Range("F476").Value = 100112009
$ws.Range("G476").Value = "Acelga"
$ws.Range("H476").Value = "Sin especificar"
$ws.Range("I476").Value = "Primera"
$ws.Range("J476").Value = 52
$ws.Range("K476").Value = 18000
$ws.Range("L476").Value = 18000
$ws.Range("M476").Value = 18000
$ws.Range("N476").Value = "`$/docena de atados"
$ws.Range("O476").Value = "Región Metropolitana"
$ws.Range("P476").Value = 6000
$ws.Range("Q476").Value = 3
$ws.Range("R476").Value = "Hortaliza"

# New row 477: Acelga, Segunda, week of 2022-02-18
$ws.Range("A477").Value = 9
$ws.Range("B477").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C477").Value = "Metropolitana"
$ws.Range("D477").Value = 44610
$ws.Range("E477").Value = 13
$ws.Range("F477").Value = 100112009
$ws.Range("G477").Value = "Acelga"
$ws.Range("H477").Value = "Sin especificar"
$ws.Range("I477").Value = "Segunda"
$ws.Range("J477").Value = 25
$ws.Range("K477").Value = 13000
$ws.Range("L477").Value = 13000
$ws.Range("M477").Value = 13000
$ws.Range("N477").Value = "`$/docena de atados"
$ws.Range("O477").Value = "Región Metropolitana"
$ws.Range("P477").Value = 4333
$ws.Range("Q477").Value = 3
$ws.Range("R477").Value = "Hortaliza"
